$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set X (pu) column (D2:D35) to the uniform value 0.002 (2E-3)
$ws.Range("D2:D35").Value = 0.002

# Update the active selection to match the saved view state
$ws.Range("I20").Select()
